$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-blank "F" column attendance figures. The I (total)
# and J (percentage) columns are formulas, so they recalculate automatically.
$ws.Range("F10").Value = 22
$ws.Range("F12").Value = 3
$ws.Range("F14").Value = 3
$ws.Range("F15").Value = 3
$ws.Range("F16").Value = 3
$ws.Range("F17").Value = 3
$ws.Range("F18").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("F20").Value = 3
$ws.Range("F21").Value = 0
$ws.Range("F22").Value = 3
$ws.Range("F23").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("F25").Value = 3

# Move the selection/scroll position down to row 13, matching the author's
# on-screen view when the file was saved.
$ws.Activate()
$ws.Range("D13:J13").Select()
